$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column E (EmiCO2) data for rows 2-11 -----------------------------
# Most of these look like numbers ("0.059", "0.078", ...) but must be stored
# as literal text, exactly like the source workbook (shared string, not a
# numeric value). Temporarily formatting the cell as Text ("@") before
# assigning the value forces Excel to keep the literal text instead of
# auto-converting it to a number; restoring the "Normal" style afterwards
# keeps the text content while dropping the temporary Text number format, so
# the cell ends up as a plain text cell with the default style.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.059"
$ws.Range("E2").Style = "Normal"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.059"
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.078"
$ws.Range("E4").Style = "Normal"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.093"
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.085"
$ws.Range("E6").Style = "Normal"

# "." is not a numeric literal, so it is stored as text automatically.
$ws.Range("E7").Value = "."

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.060"
$ws.Range("E8").Style = "Normal"

# Row 9 stays a plain number (0), matching the source data.
$ws.Range("E9").Value = 0

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.071"
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "."

# --- Column widths ----------------------------------------------------------
# Column B (Nome) and column H (Type) get explicit custom widths.
$ws.Columns(2).ColumnWidth = 13.3
$ws.Columns(8).ColumnWidth = 27.5

# --- Selection ---------------------------------------------------------------
# The active selection moves to J11.
$ws.Range("J11").Select()
